# Update to 2024 crop data:
# Adds a second round of sow-fertilization records (sow_fert_date_2 /
# sow_fert_type_npk_2 / sow_fert_type_name_2 / sow_fert_dose_2 — columns
# S:V) dated 2024-04-04, duplicating the existing first-round fertilizer
# info (2-2-2 Espoma Organic Grow!, half dose) for every crop that
# already had a first fertilization logged (rows 3-20).
#
# For the crops that had NO fertilization logged yet (tomatoes, basil,
# etc. — rows 21-33), the *first* round (sow_fert_date_1 / npk_1 /
# name_1 / dose_1 — columns O:R) is now populated with that same
# 2024-04-04 / 2-2-2 / Espoma- Organic Grow! / half record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fertDate = "2024-04-04"
$fertNpk  = "2-2-2"
$fertName = "Espoma- Organic Grow!"
$fertDose = "half"

function Set-TextValue($range, $value) {
    # Force the cell to stay a text value — otherwise Excel will parse a
    # date-shaped string like "2024-04-04" into a date serial number.
    $range.NumberFormat = "@"
    $range.Value = $value
    # Drop the number-format override again so the cell ends up on the
    # plain default style, matching the rest of the sheet's untouched
    # text cells.
    $range.Style = "Normal"
}

# Rows 3-20 already have a first fertilization entry (O:R). Duplicate it
# into the second fertilization slot (S:V) with the new date.
for ($row = 3; $row -le 20; $row++) {
    Set-TextValue $ws.Range("S$row") $fertDate
    Set-TextValue $ws.Range("T$row") $fertNpk
    Set-TextValue $ws.Range("U$row") $fertName
    Set-TextValue $ws.Range("V$row") $fertDose
}

# Rows 21-33 had no fertilization logged yet, so this becomes their
# first fertilization entry (O:R), with the new date.
for ($row = 21; $row -le 33; $row++) {
    Set-TextValue $ws.Range("O$row") $fertDate
    Set-TextValue $ws.Range("P$row") $fertNpk
    Set-TextValue $ws.Range("Q$row") $fertName
    Set-TextValue $ws.Range("R$row") $fertDose
}
